$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing value in A2: reorder "2-queque,1-torta," -> "1-torta,2-queque," ---
$ws.Range("A2").Value = "1-torta,2-queque,"

# Cells whose literal text looks like a date or a pure number need to be
# pre-formatted as Text so Excel stores them as strings instead of
# auto-converting them to a date serial / numeric value.
$textForced = @("B7","C7","F7","B8","C8","F8")
foreach ($addr in $textForced) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Add new row 7 ---
$ws.Range("A7").Value = "1-queque,"
$ws.Range("B7").Value = "10-12-2010"
$ws.Range("C7").Value = "10-12-2013"
$ws.Range("D7").Value = "test"
$ws.Range("E7").Value = "aaaa"
$ws.Range("F7").Value = "9999999"
$ws.Range("G7").Value = "Pendiente"
$ws.Range("H7").Value = 500
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0

# --- Add new row 8 ---
$ws.Range("A8").Value = "1-testito,"
$ws.Range("B8").Value = "10-10-2010"
$ws.Range("C8").Value = "10-10-2010"
$ws.Range("D8").Value = "aaaa"
$ws.Range("E8").Value = "aaa"
$ws.Range("F8").Value = "123"
$ws.Range("G8").Value = "En Proceso"
$ws.Range("H8").Value = 12
$ws.Range("I8").Value = 10
$ws.Range("J8").Value = 0

# Remove the temporary Text formatting again so the new cells don't carry a
# style reference that wasn't present in the target sheet.
foreach ($addr in $textForced) {
    $ws.Range($addr).ClearFormats()
}
